$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Labels")

# Insert a new column before the "text" column (currently H) to hold the new
# "placement" field. Inserting shifts H->I (text) and I->J (labelType), and
# the new column inherits number formats/styles from its neighbours per row.
$ws.Range("H1").EntireColumn.Insert()

# Copy the header formatting from the neighbouring "c9olumnEnd" header cell
# onto the new header cell, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "placement"

# rows -> grid now spans 3 columns instead of 2
$ws.Range("G2").Value = 3

# Move the active tab / selection to the Labels sheet, matching the diff.
$ws.Activate()
$ws.Range("H8").Select()

$wb.Save()
